$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 348.22223
$ws.Range("I4").Value = 313.42856
$ws.Range("K4").Value = 313.42856
$ws.Range("M4").Value = -199.42856
$ws.Range("H43").Value = 4346.8237
$ws.Range("I43").Value = 3366.5
$ws.Range("J43").Value = 4881.5454
$ws.Range("K43").Value = 3366.5
$ws.Range("L43").Value = 4881.5454
$ws.Range("M43").Value = -3297.5
$ws.Range("N43").Value = -5019.5454
$ws.Range("H86").Value = 2069.647
$ws.Range("I86").Value = 1717.9
$ws.Range("K86").Value = 1717.9
$ws.Range("M86").Value = -594.9000000000001
$ws.Range("H89").Value = 2069.647
$ws.Range("I89").Value = 1717.9
$ws.Range("K89").Value = 8589.5
$ws.Range("M89").Value = -2973.5
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("N93").Value = 0
$ws.Range("H98").Value = 2133.4167
$ws.Range("J98").Value = 6000
$ws.Range("L98").Value = 6000
$ws.Range("N98").Value = -8996
$ws.Range("H112").Value = 5125.5293
$ws.Range("J112").Value = 5402.4375
$ws.Range("L112").Value = 16207.3125
$ws.Range("N112").Value = -18423.3125
$ws.Range("H113").Value = 5199.4
$ws.Range("I113").Value = 3799.8
$ws.Range("J113").Value = 6599
$ws.Range("K113").Value = 3799.8
$ws.Range("L113").Value = 6599
$ws.Range("M113").Value = -545.8000000000002
$ws.Range("N113").Value = -13107
$ws.Range("H116").Value = 15802.76
$ws.Range("I116").Value = 15581.0625
$ws.Range("K116").Value = 15581.0625
$ws.Range("M116").Value = -12139.0625
$ws.Range("H122").Value = 2133.4167
$ws.Range("J122").Value = 6000
$ws.Range("L122").Value = 18000
$ws.Range("N122").Value = -22900
$ws.Range("H132").Value = 37271.35
$ws.Range("I132").Value = 41035.168
$ws.Range("K132").Value = 123105.504
$ws.Range("M132").Value = -120575.504
$ws.Range("H137").Value = 12458.034
$ws.Range("I137").Value = 15585.909
$ws.Range("K137").Value = 46757.727
$ws.Range("M137").Value = -44207.727
$ws.Range("H138").Value = 28122.281
$ws.Range("I138").Value = 1800.5714
$ws.Range("J138").Value = 58830.945
$ws.Range("K138").Value = 5401.7142
$ws.Range("L138").Value = 176492.835
$ws.Range("M138").Value = -261.7142000000003
$ws.Range("N138").Value = -186772.835
$ws.Range("H139").Value = 99999
$ws.Range("J139").Value = 99999
$ws.Range("L139").Value = 99999
$ws.Range("N139").Value = -110279
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 38615.83
$ws.Range("I32").Value = 42686.883
$ws.Range("J32").Value = 3333.3333
$ws.Range("K32").Value = 42686.883
$ws.Range("L32").Value = 3333.3333
$ws.Range("M32").Value = -42399.883
$ws.Range("N32").Value = -3907.3333
$ws.Range("H61").Value = 4239
$ws.Range("I61").Value = 1068.04
$ws.Range("K61").Value = 1068.04
$ws.Range("M61").Value = -856.04
$ws.Range("H132").Value = 2661.8823
$ws.Range("I132").Value = 2154.5557
$ws.Range("K132").Value = 6463.6671
$ws.Range("M132").Value = -3933.6671
$ws.Range("H136").Value = 4239
$ws.Range("I136").Value = 1068.04
$ws.Range("K136").Value = 3204.12
$ws.Range("M136").Value = -654.1199999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2683.5
$ws.Range("I64").Value = 2266.5
$ws.Range("K64").Value = 2266.5
$ws.Range("M64").Value = -2041.5
$ws.Range("H67").Value = 2683.5
$ws.Range("I67").Value = 2266.5
$ws.Range("K67").Value = 2266.5
$ws.Range("M67").Value = -1486.5
$ws.Range("H134").Value = 2873.7827
$ws.Range("I134").Value = 2189.5625
$ws.Range("K134").Value = 6568.6875
$ws.Range("M134").Value = -4033.6875
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 11669
$ws.Range("I17").Value = 11669
$ws.Range("K17").Value = 11669
$ws.Range("M17").Value = -11495
$ws.Range("H99").Value = 6273.8125
$ws.Range("I99").Value = 4898.2
$ws.Range("K99").Value = 4898.2
$ws.Range("M99").Value = -3400.2
$ws.Range("H107").Value = 724.2
$ws.Range("I107").Value = 666.6875
$ws.Range("K107").Value = 666.6875
$ws.Range("M107").Value = 1253.3125
$ws.Range("H126").Value = 6273.8125
$ws.Range("I126").Value = 4898.2
$ws.Range("K126").Value = 14694.6
$ws.Range("M126").Value = -12224.6
$ws.Range("H132").Value = 57144.445
$ws.Range("I132").Value = 77807
$ws.Range("J132").Value = 3421.8
$ws.Range("K132").Value = 233421
$ws.Range("L132").Value = 10265.4
$ws.Range("M132").Value = -230891
$ws.Range("N132").Value = -15325.4
$ws.Range("H134").Value = 2141.25
$ws.Range("I134").Value = 1823.6111
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 5470.8333
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -2935.8333
$ws.Range("N134").Value = -20070
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 41972.555
$ws.Range("J37").Value = 41972.555
$ws.Range("L37").Value = 125917.665
$ws.Range("N37").Value = -126141.665
$ws.Range("H69").Value = 4408.706
$ws.Range("I69").Value = 2249.5
$ws.Range("J69").Value = 4696.6
$ws.Range("K69").Value = 6748.5
$ws.Range("L69").Value = 14089.8
$ws.Range("M69").Value = -5937.5
$ws.Range("N69").Value = -15711.8
$ws.Range("H72").Value = 4408.706
$ws.Range("I72").Value = 2249.5
$ws.Range("J72").Value = 4696.6
$ws.Range("K72").Value = 20245.5
$ws.Range("L72").Value = 42269.4
$ws.Range("M72").Value = -16189.5
$ws.Range("N72").Value = -50381.4
$ws.Range("H131").Value = 120554.36
$ws.Range("I131").Value = 428894.3
$ws.Range("J131").Value = 1962.0769
$ws.Range("K131").Value = 1286682.9
$ws.Range("L131").Value = 5886.2307
$ws.Range("M131").Value = -1281642.9
$ws.Range("N131").Value = -15966.2307
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = 0
$ws.Range("H132").Value = 1990.2122
$ws.Range("I132").Value = 1689.2667
$ws.Range("J132").Value = 4999.6665
$ws.Range("K132").Value = 5067.800099999999
$ws.Range("L132").Value = 14998.9995
$ws.Range("M132").Value = -2537.800099999999
$ws.Range("N132").Value = -20058.9995
$ws.Range("H140").Value = 84800
$ws.Range("I140").Value = 60000
$ws.Range("J140").Value = 91000
$ws.Range("K140").Value = 60000
$ws.Range("L140").Value = 91000
$ws.Range("M140").Value = -54820
$ws.Range("N140").Value = -101360
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2724.0527
$ws.Range("I7").Value = 2724.0527
$ws.Range("K7").Value = 2724.0527
$ws.Range("M7").Value = -2612.0527
$ws.Range("H16").Value = 843.1875
$ws.Range("I16").Value = 1117.1818
$ws.Range("J16").Value = 240.4
$ws.Range("K16").Value = 1117.1818
$ws.Range("L16").Value = 240.4
$ws.Range("M16").Value = -947.1818000000001
$ws.Range("N16").Value = -580.4
$ws.Range("H46").Value = 4715.4287
$ws.Range("I46").Value = 970
$ws.Range("J46").Value = 6796.222
$ws.Range("K46").Value = 970
$ws.Range("L46").Value = 6796.222
$ws.Range("M46").Value = -782
$ws.Range("N46").Value = -7172.222
$ws.Range("H126").Value = 2724.0527
$ws.Range("I126").Value = 2724.0527
$ws.Range("K126").Value = 8172.158100000001
$ws.Range("M126").Value = -5702.158100000001
$ws.Range("H132").Value = 2118.8823
$ws.Range("I132").Value = 932.8570999999999
$ws.Range("J132").Value = 2949.1
$ws.Range("K132").Value = 2798.5713
$ws.Range("L132").Value = 8847.299999999999
$ws.Range("M132").Value = -268.5712999999996
$ws.Range("N132").Value = -13907.3
$ws.Range("H136").Value = 4410.3706
$ws.Range("I136").Value = 4006.6667
$ws.Range("J136").Value = 5823.3335
$ws.Range("K136").Value = 12020.0001
$ws.Range("L136").Value = 17470.0005
$ws.Range("M136").Value = -9470.000100000001
$ws.Range("N136").Value = -22570.0005
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 999
$ws.Range("I2").Value = 999
$ws.Range("K2").Value = 999
$ws.Range("M2").Value = -887
$ws.Range("H132").Value = 1698.6296
$ws.Range("I132").Value = 1272.3684
$ws.Range("K132").Value = 3817.1052
$ws.Range("M132").Value = -1287.1052
$ws.Range("H136").Value = 22605.96
$ws.Range("I136").Value = 25784.363
$ws.Range("K136").Value = 77353.08900000001
$ws.Range("M136").Value = -74803.08900000001
